$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 543, shifting the existing 543:574 block down to 546:577.
$ws.Range("A543:A545").EntireRow.Insert()

# New weekly record (Fecha 44516) split across the three usual "Calidad" tiers,
# following the same layout as every other row in the table.
$data = @(
    @(543, "Primera", 57000, 70, 80, 74),
    @(544, "Segunda", 49000, 55, 60, 57),
    @(545, "Tercera", 15000, 35, 35, 35)
)

foreach ($row in $data) {
    $r = $row[0]
    $calidad = $row[1]
    $volumen = $row[2]
    $pmin = $row[3]
    $pmax = $row[4]
    $pprom = $row[5]

    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = 44516
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100114014
    $ws.Cells.Item($r, 7).Value = "Betarraga"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $pmin
    $ws.Cells.Item($r, 12).Value = $pmax
    $ws.Cells.Item($r, 13).Value = $pprom
    $ws.Cells.Item($r, 14).Value = "$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
